$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row text / value fix ---
$ws.Range("D1").Value2 = "Supplier"

# --- Cell text corrections throughout the BOM ---
$ws.Range("A2").Value2  = "Photoreciever"
$ws.Range("A12").Value2 = "M6 - M3 thread adapter"
$ws.Range("C13").Value2 = "CL3/M Compact Variable Height Clamp, M6 Tapped"
$ws.Range("A16").Value2 = "M6 screw 12mm "
$ws.Range("C16").Value2 = "SH6MS12"
$ws.Range("A19").Value2 = "M3 screw 10mm"
$ws.Range("A20").Value2 = "M3 spacer 3mm"
$ws.Range("A21").Value2 = "M6 x 45mm"

# --- Make header row bold ---
$ws.Range("A1:D1").Font.Bold = $true

# --- Column width tweaks ---
# (ColumnWidth is quantized internally to 1/6-character steps by the COM
# layer, so these inputs are chosen to land on the stored width closest to
# the target: A -> 51.333333..., C -> 52.5)
$ws.Columns.Item(1).ColumnWidth = 50.5
$ws.Columns.Item(3).ColumnWidth = 51.6665

# --- Update active cell selection ---
$ws.Range("D1").Select()
